# Add a new "2020" data column (column N) to the Sheet1 table, mirroring
# the formatting already used by the other year columns, and move the
# active selection as recorded by Excel when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (row 4): style copied from the thick-bordered
# header cell used by the rest of the year row (D4:M4), then rendered bold.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# --- Row 3 (separator row under the title): extend the thin bottom
# border one more column so it lines up under the new header cell.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# --- Data rows 5-9: plain values with the same font as the rest of the
# data columns, no outer border.
$ws.Range("L5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 1713

$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 1

$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 379

$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 180

# --- Row 10 (bottom, thick-bottom-bordered totals row)
$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 798

# --- Restore the selection Excel had recorded for this sheet the last
# time the file was saved.
$ws.Range("L22").Select()
